$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Corrección error reporte ventas: el pedido de Isavo Castro (fila 7)
# figuraba como "Finalizado" y debía estar como "Retirado".
$ws.Range("G7").Value = "Retirado"
